# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/mora table (rows 16-36) is re-sorted: originally grouped by
# worker (Balvina, Joaquin, Daivis) with periods descending; now grouped by
# period (ascending: 1911, 1912, 2001..2005) with the same three workers
# (in the same relative order) repeated under each period. No new people or
# amounts are introduced - it's the same 21 rows, re-ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45555741"
$ws.Range("D16").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E16").Value = "1911"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1083452289"
$ws.Range("D17").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E17").Value = "1911"
$ws.Range("F17").Value = 60000
$ws.Range("G17").Value = 1500000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "52708490"
$ws.Range("D18").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E18").Value = "1911"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 828116

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45555741"
$ws.Range("D19").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E19").Value = "1912"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1083452289"
$ws.Range("D20").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E20").Value = "1912"
$ws.Range("F20").Value = 60000
$ws.Range("G20").Value = 1500000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "52708490"
$ws.Range("D21").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E21").Value = "1912"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 828116

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45555741"
$ws.Range("D22").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E22").Value = "2001"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1083452289"
$ws.Range("D23").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E23").Value = "2001"
$ws.Range("F23").Value = 60000
$ws.Range("G23").Value = 1500000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "52708490"
$ws.Range("D24").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E24").Value = "2001"
$ws.Range("F24").Value = 33125
$ws.Range("G24").Value = 828116

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "45555741"
$ws.Range("D25").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E25").Value = "2002"
$ws.Range("F25").Value = 33125
$ws.Range("G25").Value = 828116

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1083452289"
$ws.Range("D26").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E26").Value = "2002"
$ws.Range("F26").Value = 60000
$ws.Range("G26").Value = 1500000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "52708490"
$ws.Range("D27").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E27").Value = "2002"
$ws.Range("F27").Value = 33125
$ws.Range("G27").Value = 828116

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "45555741"
$ws.Range("D28").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E28").Value = "2003"
$ws.Range("F28").Value = 33125
$ws.Range("G28").Value = 828116

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1083452289"
$ws.Range("D29").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E29").Value = "2003"
$ws.Range("F29").Value = 60000
$ws.Range("G29").Value = 1500000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "52708490"
$ws.Range("D30").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E30").Value = "2003"
$ws.Range("F30").Value = 33125
$ws.Range("G30").Value = 828116

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "45555741"
$ws.Range("D31").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E31").Value = "2004"
$ws.Range("F31").Value = 33125
$ws.Range("G31").Value = 828116

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1083452289"
$ws.Range("D32").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E32").Value = "2004"
$ws.Range("F32").Value = 60000
$ws.Range("G32").Value = 1500000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "52708490"
$ws.Range("D33").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E33").Value = "2004"
$ws.Range("F33").Value = 33125
$ws.Range("G33").Value = 828116

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "45555741"
$ws.Range("D34").Value = "BALVINA MEZA BALLESTEROS"
$ws.Range("E34").Value = "2005"
$ws.Range("F34").Value = 20979
$ws.Range("G34").Value = 828116

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1083452289"
$ws.Range("D35").Value = "JOAQUIN ANGEL MENDOZA SILVA"
$ws.Range("E35").Value = "2005"
$ws.Range("F35").Value = 38000
$ws.Range("G35").Value = 1500000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "52708490"
$ws.Range("D36").Value = "DAIVIS ANGELICA ROJAS PAREDES"
$ws.Range("E36").Value = "2005"
$ws.Range("F36").Value = 20979
$ws.Range("G36").Value = 828116
